$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.929.97'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +5.94%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.230.88'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.93%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '231.63'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.94%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '61.02'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -2.91%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '58.93'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.98%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +4.86%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.23%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.562.04'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.90%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.63'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -1.48%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.75'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.25%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.797'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.08%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.58'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.08%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.250.46'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.67%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '41.776.69'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +5.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.37'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.02%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.77%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.03'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.44%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '249.98'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +9.75%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.39'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.69%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.29%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.70'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.143'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.98%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '167.00'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.97%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.95'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.41'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.54%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.121'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.37%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.95'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +5.34%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.63'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +3.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0630'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.73%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.64'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -4.68%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.70%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.36'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.53%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.000256'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +29.74%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0241'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +6.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.81'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.25%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +8.49%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0979'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +7.16%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '98.88'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -3.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.22'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.472.97'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.66%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '16.51'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -7.02%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.11%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.21%  '
